# Generate Report for Handback
#
# Updates the localization-status report:
#  - Row 3 (fe2f219f-...) status moves from "Ready for handoff" to
#    "Handback transform failed" (shared across Overview/zh-cn/de-de sheets).
#  - Adds an "Error Detail" note (column L) on row 3 of the zh-cn and
#    de-de sheets explaining the handback file name mismatch.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update the status text for the fe2f219f row (row 3) everywhere it appears.
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Add the Error Detail message (column L) for row 3 on the zh-cn sheet.
$wsZhCn.Range("L3").Value = "Handback file name: 3u3yzoal.3me is different with handoff file name: fe2f219f-a663-42ce-8593-0e58e8d2024d.0db1f6a176ef375b40b7772d3b3b57539010b350.zh-cn."

# Add the Error Detail message (column L) for row 3 on the de-de sheet.
$wsDeDe.Range("L3").Value = "Handback file name: 3u3yzoal.3me is different with handoff file name: fe2f219f-a663-42ce-8593-0e58e8d2024d.0db1f6a176ef375b40b7772d3b3b57539010b350.de-de."
